# Auto-generated PowerShell Word COM-interop edit script
$d = $word.ActiveDocument

# Step 1: locate the final (empty, centered) paragraph that will become the first remediation bullet
$count = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($count)

# Step 2: give it simple text + apply Word's default numbered-list formatting.
# This creates a brand new abstractNum/num pair (numId referenced below) in numbering.xml,
# mirroring how the author generated a fresh numbered list for the remediation bullets.
$target.Range.Text = "x"
$target.Range.ListFormat.ApplyNumberDefault()

# Step 3: shape the new list template's levels so they follow the same decimal / lowerLetter / lowerRoman
# rotation used by the document's existing multilevel list (abstractNumId 0).
$newTemplate = $target.Range.ListFormat.ListTemplate
$numberStylePattern = @(0,4,2,0,4,2,0,4,2)
for ($lvlIdx = 1; $lvlIdx -le 9; $lvlIdx++) {
    $lvl = $newTemplate.ListLevels.Item($lvlIdx)
    $lvl.NumberStyle = $numberStylePattern[$lvlIdx - 1]
}

# Step 4: replace the (now list-formatted) paragraph's content with the first remediation bullet,
# using InsertXML so the exact run/paragraph formatting from the source matches precisely.
$targetRange = $d.Paragraphs.Item($count).Range
$p1 = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>Consistent Request Header Parsing:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Both the front-end and back-end servers should consistently parse and process HTTP request headers. In the event of duplicate headers, either reject the request or combine them as per the HTTP specification.</w:t></w:r></w:p></pkg:xmlData>
'@
$targetRange.InsertXML($p1)

# Step 5: append the remaining five remediation bullets as new paragraphs, each inserted
# immediately after the current last paragraph and then populated via InsertXML.
$prevCount = $d.Paragraphs.Count
$prevPara = $d.Paragraphs.Item($prevCount)
$prevPara.Range.InsertParagraphAfter()
$newCount = $d.Paragraphs.Count
$newRange = $d.Paragraphs.Item($newCount).Range
$p2 = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>Prohibit Unknown Transfer Encodings:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> If the server encounters an unknown transfer encoding like "cow" in the above example, it should reject the request. Servers should only process known and valid encodings.</w:t></w:r></w:p></pkg:xmlData>
'@
$newRange.InsertXML($p2)

$prevCount = $d.Paragraphs.Count
$prevPara = $d.Paragraphs.Item($prevCount)
$prevPara.Range.InsertParagraphAfter()
$newCount = $d.Paragraphs.Count
$newRange = $d.Paragraphs.Item($newCount).Range
$p3 = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>Implement a Strict Allowlist:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Rather than blocking known harmful headers or encodings, implement an allowlist approach. Only explicitly allowed headers and encodings get processed, and everything else gets rejected.</w:t></w:r></w:p></pkg:xmlData>
'@
$newRange.InsertXML($p3)

$prevCount = $d.Paragraphs.Count
$prevPara = $d.Paragraphs.Item($prevCount)
$prevPara.Range.InsertParagraphAfter()
$newCount = $d.Paragraphs.Count
$newRange = $d.Paragraphs.Item($newCount).Range
$p4 = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>Use Web Application Firewalls (WAFs):</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> WAFs can detect anomalies in the HTTP requests. Configuring a WAF to recognize and block such anomalies can prevent HTTP Request Smuggling attacks.</w:t></w:r></w:p></pkg:xmlData>
'@
$newRange.InsertXML($p4)

$prevCount = $d.Paragraphs.Count
$prevPara = $d.Paragraphs.Item($prevCount)
$prevPara.Range.InsertParagraphAfter()
$newCount = $d.Paragraphs.Count
$newRange = $d.Paragraphs.Item($newCount).Range
$p5 = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>Regularly Update and Patch Servers:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Both front-end and back-end servers should be updated and patched regularly. Many HTTP Request Smuggling vulnerabilities arise from outdated server software.</w:t></w:r></w:p></pkg:xmlData>
'@
$newRange.InsertXML($p5)

Write-Host "Remediation section inserted. Paragraph count now: $($d.Paragraphs.Count)"
